$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.497.51"
$ws.Range("E2").Value = "  -2.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.591.20"
$ws.Range("E3").Value = "  -2.18%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.82"
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.46"
$ws.Range("E6").Value = "  -2.58%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  -2.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.599.33"
$ws.Range("E9").Value = "  -2.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.68"
$ws.Range("E10").Value = "  -2.50%  "
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("E12").Value = "  +10.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.360"
$ws.Range("E13").Value = "  +4.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.047.67"
$ws.Range("E14").Value = "  -2.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.40"
$ws.Range("E15").Value = "  +6.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "59.428.31"
$ws.Range("E16").Value = "  -2.41%  "
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.591.06"
$ws.Range("E18").Value = "  -2.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.60"
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "338.63"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.44"
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.46"
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.91"
$ws.Range("E24").Value = "  -4.09%  "
$ws.Range("E25").Value = "  +6.23%  "
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("E27").Value = "  -2.79%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0781"
$ws.Range("E29").Value = "  -0.68%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("E32").Value = "  -2.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "158.21"
$ws.Range("E33").Value = "  +2.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.08"
$ws.Range("E34").Value = "  -0.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.06"
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.887"
$ws.Range("E37").Value = "  -2.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.876"
$ws.Range("E38").Value = "  -3.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.45"
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("E40").Value = "  -1.99%  "
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "294.17"
$ws.Range("E42").Value = "  -3.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "135.74"
$ws.Range("E43").Value = "  +5.19%  "
$ws.Range("E44").Value = "  +0.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0978"
$ws.Range("E45").Value = "  -0.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.598"
$ws.Range("E46").Value = "  -1.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0537"
$ws.Range("E47").Value = "  -2.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.63"
$ws.Range("E48").Value = "  -0.21%  "
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.87"
$ws.Range("E50").Value = "  -0.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.956.05"
$ws.Range("E51").Value = "  -0.67%  "
